# Applies the commit's two changes:
#   1. The table on slide 5 switches from the "Table_0" style
#      ({6DC17D52-40D0-447D-A1DC-3D0FA883AC8B}) to the
#      {8B6219EC-1076-4DB5-80DC-FBE9060A74B2} style.
#   2. The deck's theme colour scheme (which currently carries the
#      "Integral" / "Red Violet" palette) is swapped for the classic
#      "Office Theme" palette, mirroring the theme1.xml <-> theme2.xml
#      content swap in the target OOXML.

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Table style id on the slide-5 table
# -----------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{8B6219EC-1076-4DB5-80DC-FBE9060A74B2}")

# -----------------------------------------------------------------
# 2) Theme colour swap (Integral/Red Violet -> Office Theme)
# -----------------------------------------------------------------
function ConvertTo-VbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the standard theme colour scheme layout:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $colorScheme.Item($i + 1).RGB = ConvertTo-VbaRgb $officeThemeColors[$i]
}
